# Auto-generated edit script applying scheduled-runner value updates
# to the Leve profit-tracking sheets (currentAveragePrice* / LevePrice* / LeveProfit* columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53: No Accounting for Waste | Enchanted Electrum Ink
$ws.Range("H53").Value2 = 349.5
$ws.Range("I53").Value2 = 199
$ws.Range("J53").Value2 = 500
$ws.Range("K53").Value2 = 199
$ws.Range("L53").Value2 = 500
$ws.Range("M53").Value2 = 438
$ws.Range("N53").Value2 = -1774

# Row 69: Steeling the Knife, Steeling the Mind | Grade 1 Mind Dissolvent
$ws.Range("H69").Value2 = 19675.375
$ws.Range("J69").Value2 = 23166.666
$ws.Range("L69").Value2 = 69499.99800000001
$ws.Range("N69").Value2 = -71247.99800000001

# Row 72: Surgical Substitution (L) | Grade 1 Mind Dissolvent
$ws.Range("H72").Value2 = 19675.375
$ws.Range("J72").Value2 = 23166.666
$ws.Range("L72").Value2 = 208499.994
$ws.Range("N72").Value2 = -217235.994

# Row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws.Range("H76").Value2 = 3772.0908
$ws.Range("I76").Value2 = 2998.625
$ws.Range("K76").Value2 = 2998.625
$ws.Range("M76").Value2 = -2683.625

# Row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws.Range("H79").Value2 = 3772.0908
$ws.Range("I79").Value2 = 2998.625
$ws.Range("K79").Value2 = 2998.625
$ws.Range("M79").Value2 = -1906.625

# Row 113: Amaro Kart | Starch Glue
$ws.Range("H113").Value2 = 3294
$ws.Range("I113").Value2 = 3492.5
$ws.Range("K113").Value2 = 3492.5
$ws.Range("M113").Value2 = -238.5

# Row 140: Tome for Tradition | Book of Ra'Kaznar
$ws.Range("H140").Value2 = 98988
$ws.Range("J140").Value2 = 98988
$ws.Range("L140").Value2 = 98988
$ws.Range("N140").Value2 = -109348

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value2 = 1416.625
$ws.Range("I32").Value2 = 1397.1
$ws.Range("K32").Value2 = 1397.1
$ws.Range("M32").Value2 = -1110.1

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value2 = 8338305
$ws.Range("I132").Value2 = 10003763
$ws.Range("K132").Value2 = 30011289
$ws.Range("M132").Value2 = -30008759

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value2 = 2353.3572
$ws.Range("I86").Value2 = 2380.6924
$ws.Range("K86").Value2 = 2380.6924
$ws.Range("M86").Value2 = -1257.6924

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value2 = 2353.3572
$ws.Range("I89").Value2 = 2380.6924
$ws.Range("K89").Value2 = 11903.462
$ws.Range("M89").Value2 = -6287.462

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value2 = 999999
$ws.Range("I107").Value2 = 0
$ws.Range("J107").Value2 = 999999
$ws.Range("K107").Value2 = 0
$ws.Range("L107").Value2 = 999999
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value2 = -1003839

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value2 = 6665.2583
$ws.Range("I31").Value2 = 4749.5454
$ws.Range("J31").Value2 = 7718.9
$ws.Range("K31").Value2 = 4749.5454
$ws.Range("L31").Value2 = 7718.9
$ws.Range("M31").Value2 = -4454.5454
$ws.Range("N31").Value2 = -8308.9

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value2 = 6665.2583
$ws.Range("I34").Value2 = 4749.5454
$ws.Range("J34").Value2 = 7718.9
$ws.Range("K34").Value2 = 4749.5454
$ws.Range("L34").Value2 = 7718.9
$ws.Range("M34").Value2 = -4547.5454
$ws.Range("N34").Value2 = -8122.9

# Row 35: Storm of Swords | Elm Macuahuitl
$ws.Range("H35").Value2 = 2216.7
$ws.Range("I35").Value2 = 1019
$ws.Range("K35").Value2 = 1019
$ws.Range("M35").Value2 = -725

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value2 = 25005564
$ws.Range("J58").Value2 = 169
$ws.Range("L58").Value2 = 169
$ws.Range("N58").Value2 = -575

# Row 98: Pinewheel | Pine Spinning Wheel
$ws.Range("H98").Value2 = 103333
$ws.Range("J98").Value2 = 103333
$ws.Range("L98").Value2 = 103333
$ws.Range("N98").Value2 = -107825

# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value2 = 2649.7827
$ws.Range("I122").Value2 = 2952.4211
$ws.Range("J122").Value2 = 1212.25
$ws.Range("K122").Value2 = 8857.263300000001
$ws.Range("L122").Value2 = 3636.75
$ws.Range("M122").Value2 = -6407.263300000001
$ws.Range("N122").Value2 = -8536.75

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value2 = 19609676
$ws.Range("I132").Value2 = 24391728
$ws.Range("K132").Value2 = 73175184
$ws.Range("M132").Value2 = -73172654

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value2 = 6946048
$ws.Range("I134").Value2 = 7577318.5
$ws.Range("J134").Value2 = 2074.3333
$ws.Range("K134").Value2 = 22731955.5
$ws.Range("L134").Value2 = 6222.999899999999
$ws.Range("M134").Value2 = -22729420.5
$ws.Range("N134").Value2 = -11292.9999

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value2 = 25005564
$ws.Range("J136").Value2 = 169
$ws.Range("L136").Value2 = 507
$ws.Range("N136").Value2 = -5607

$ws = $wb.Worksheets.Item("CUL")
# Row 23: Sweet Smell of Success | Lavender Oil
$ws.Range("H23").Value2 = 754.0526
$ws.Range("J23").Value2 = 857.4167
$ws.Range("L23").Value2 = 2572.2501
$ws.Range("N23").Value2 = -3042.2501

# Row 60: Drinking to Your Health | Mulled Tea
$ws.Range("H60").Value2 = 4414.905
$ws.Range("I60").Value2 = 1046.6364
$ws.Range("J60").Value2 = 8120
$ws.Range("K60").Value2 = 3139.9092
$ws.Range("L60").Value2 = 24360
$ws.Range("M60").Value2 = -2888.9092
$ws.Range("N60").Value2 = -24862

# Row 68: Such a Butter Face | Fermented Butter
$ws.Range("H68").Value2 = 1014.1667
$ws.Range("J68").Value2 = 1014.1667
$ws.Range("L68").Value2 = 3042.5001
$ws.Range("N68").Value2 = -4664.5001

# Row 71: No Margarine of Error (L) | Fermented Butter
$ws.Range("H71").Value2 = 1014.1667
$ws.Range("J71").Value2 = 1014.1667
$ws.Range("L71").Value2 = 9127.5003
$ws.Range("N71").Value2 = -17239.5003

# Row 134: Don't Knock It Till You've Tried It | Mezcal-marinated Swampmonk
$ws.Range("H134").Value2 = 1304.0834
$ws.Range("I134").Value2 = 1304.0834
$ws.Range("K134").Value2 = 3912.2502
$ws.Range("M134").Value2 = 1157.7498

# Row 136: Simple Is Hardest | Spaghetti al Olio e Peperoncino
$ws.Range("H136").Value2 = 2049.8333
$ws.Range("I136").Value2 = 1559.8
$ws.Range("K136").Value2 = 4679.4
$ws.Range("M136").Value2 = 420.6000000000004

# Row 139: Najoothie | Wild Banana Blend
$ws.Range("H139").Value2 = 1209.2307
$ws.Range("I139").Value2 = 1060
$ws.Range("K139").Value2 = 3180
$ws.Range("M139").Value2 = 1960

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers | Copper Ingot
$ws.Range("H2").Value2 = 506.58334
$ws.Range("I2").Value2 = 702.0625
$ws.Range("J2").Value2 = 115.625
$ws.Range("K2").Value2 = 702.0625
$ws.Range("L2").Value2 = 115.625
$ws.Range("M2").Value2 = -589.0625
$ws.Range("N2").Value2 = -341.625

# Row 33: Thaumaturge Is Magic | Fluorite Ring
$ws.Range("H33").Value2 = 0
$ws.Range("I33").Value2 = 0
$ws.Range("K33").Value2 = 0
$ws.Range("M33").ClearContents()

# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value2 = 3109.7778
$ws.Range("I80").Value2 = 2747.25
$ws.Range("J80").Value2 = 3399.8
$ws.Range("K80").Value2 = 2747.25
$ws.Range("L80").Value2 = 3399.8
$ws.Range("M80").Value2 = -1749.25
$ws.Range("N80").Value2 = -5395.8

# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value2 = 3109.7778
$ws.Range("I83").Value2 = 2747.25
$ws.Range("J83").Value2 = 3399.8
$ws.Range("K83").Value2 = 13736.25
$ws.Range("L83").Value2 = 16999
$ws.Range("M83").Value2 = -8744.25
$ws.Range("N83").Value2 = -26983

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value2 = 4535.731
$ws.Range("I102").Value2 = 2455.375
$ws.Range("K102").Value2 = 2455.375
$ws.Range("M102").Value2 = -833.375

# Row 134: Guaranteed Gem | Ihuykanite
$ws.Range("H134").Value2 = 86666.336
$ws.Range("J134").Value2 = 86666.336
$ws.Range("L134").Value2 = 259999.008
$ws.Range("N134").Value2 = -265069.008

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value2 = 4004.4
$ws.Range("I7").Value2 = 4004.4
$ws.Range("K7").Value2 = 4004.4
$ws.Range("M7").Value2 = -3892.4

# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value2 = 6512.25
$ws.Range("I40").Value2 = 6442.5713
$ws.Range("K40").Value2 = 6442.5713
$ws.Range("M40").Value2 = -6306.5713

# Row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value2 = 1082.1666
$ws.Range("J46").Value2 = 997.5
$ws.Range("L46").Value2 = 997.5
$ws.Range("N46").Value2 = -1373.5

# Row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws.Range("H55").Value2 = 398.77777
$ws.Range("I55").Value2 = 412.7143
$ws.Range("K55").Value2 = 412.7143
$ws.Range("M55").Value2 = -239.7143

# Row 61: Spelling Me Softly | Raptor Leather
$ws.Range("H61").Value2 = 2487.5625
$ws.Range("I61").Value2 = 2320.1333
$ws.Range("J61").Value2 = 4999
$ws.Range("K61").Value2 = 2320.1333
$ws.Range("L61").Value2 = 4999
$ws.Range("M61").Value2 = -2118.1333
$ws.Range("N61").Value2 = -5403

# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Range("H113").Value2 = 2487.5625
$ws.Range("I113").Value2 = 2320.1333
$ws.Range("J113").Value2 = 4999
$ws.Range("K113").Value2 = 2320.1333
$ws.Range("L113").Value2 = 4999
$ws.Range("M113").Value2 = -150.1333
$ws.Range("N113").Value2 = -9339

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value2 = 3238.0645
$ws.Range("J122").Value2 = 8729.5
$ws.Range("L122").Value2 = 26188.5
$ws.Range("N122").Value2 = -31088.5

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value2 = 4004.4
$ws.Range("I126").Value2 = 4004.4
$ws.Range("K126").Value2 = 12013.2
$ws.Range("M126").Value2 = -9543.200000000001

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value2 = 5106169
$ws.Range("I132").Value2 = 5559741
$ws.Range("J132").Value2 = 3487.5
$ws.Range("K132").Value2 = 16679223
$ws.Range("L132").Value2 = 10462.5
$ws.Range("M132").Value2 = -16676693
$ws.Range("N132").Value2 = -15522.5

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value2 = 2159.8696
$ws.Range("I136").Value2 = 2130.4375
$ws.Range("J136").Value2 = 2227.1428
$ws.Range("K136").Value2 = 6391.3125
$ws.Range("L136").Value2 = 6681.428400000001
$ws.Range("M136").Value2 = -3841.3125
$ws.Range("N136").Value2 = -11781.4284

# Row 140: Worqor Zormor or Bust | Gargantuaskin Shoes of Healing
$ws.Range("H140").Value2 = 89765.60000000001
$ws.Range("J140").Value2 = 89765.60000000001
$ws.Range("L140").Value2 = 89765.60000000001
$ws.Range("N140").Value2 = -100125.6

$ws = $wb.Worksheets.Item("WVR")
# Row 4: Not Cool Enough | Hempen Undershirt
$ws.Range("H4").Value2 = 42939490
$ws.Range("I4").Value2 = 143800
$ws.Range("J4").Value2 = 100000400
$ws.Range("K4").Value2 = 143800
$ws.Range("L4").Value2 = 100000400
$ws.Range("M4").Value2 = -143687
$ws.Range("N4").Value2 = -100000626

# Row 107: Flax Wax | Bright Linen Yarn
$ws.Range("H107").Value2 = 3013.7144
$ws.Range("I107").Value2 = 2219.4
$ws.Range("J107").Value2 = 4999.5
$ws.Range("K107").Value2 = 6658.200000000001
$ws.Range("L107").Value2 = 14998.5
$ws.Range("M107").Value2 = -4738.200000000001
$ws.Range("N107").Value2 = -18838.5

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value2 = 15628244
$ws.Range("I132").Value2 = 22730556
$ws.Range("J132").Value2 = 3160.1
$ws.Range("K132").Value2 = 68191668
$ws.Range("L132").Value2 = 9480.299999999999
$ws.Range("M132").Value2 = -68189138
$ws.Range("N132").Value2 = -14540.3

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value2 = 8335259.5
$ws.Range("I136").Value2 = 8622510
$ws.Range("K136").Value2 = 25867530
$ws.Range("M136").Value2 = -25864980
